$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches inlineStr cells in source)
$textCells = @("D5", "D6", "D7", "D10", "D12", "D13", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D30", "D32", "D33", "D34", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values from the crypto data refresh
$ws.Range("D2").Value = "57.998.85"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.475.18"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "518.34"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "130.97"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("D9").Value = "2.506.90"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").Value = "0.0968"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "5.15"
$ws.Range("E12").Value = "  -2.85%  "
$ws.Range("D13").Value = "0.330"
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("D14").Value = "2.935.53"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "58.135.70"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "21.95"
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("D18").Value = "2.498.64"
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("D19").Value = "10.67"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "320.08"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "4.15"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.95"
$ws.Range("E23").Value = "  +4.44%  "
$ws.Range("D24").Value = "64.29"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").Value = "0.399"
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("D26").Value = "0.994"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").Value = "7.29"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "0.0₃0743"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").Value = "167.65"
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "6.21"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.17"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").Value = "17.98"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("D38").Value = "3.90"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").Value = "36.74"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").Value = "1.44"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("D41").Value = "0.765"
$ws.Range("E41").Value = "  -2.81%  "
$ws.Range("D42").Value = "273.49"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").Value = "5.05"
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("D44").Value = "3.41"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").Value = "0.593"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").Value = "0.0918"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("D47").Value = "120.31"
$ws.Range("E47").Value = "  -4.94%  "
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("D49").Value = "17.67"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("D51").Value = "16.69"
$ws.Range("E51").Value = "  -0.92%  "
